$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New way-point rows X1-X5 (rows 34-38), mirroring the existing table layout
$data = @(
    @("X1", 5.17, 2.38, "6;14"),
    @("X2", 6.35, 2.38, "8;16"),
    @("X3", 7.72, 2.38, "12;18"),
    @("X4", 5.68, 5.55, "23;29"),
    @("X5", 7.04, 5.23, "26;31")
)

$row = 34
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Formula = '=B' + $row + '+($G$2/2)'
    $ws.Cells.Item($row, 5).Formula = '=C' + $row + '+($G$2/2)'
    $ws.Cells.Item($row, 6).Value = $item[3]
    $row++
}

$ws.Range("F38").Select()
